$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strategy Worksheet")

# --- Header notes (row 1): audience + caveat ---
$ws.Range("D1").Value = "Audience = direct line manager"
$ws.Range("E1").Value = "Some factors will be skipped as out of influence of line manager"

# --- Row 4: note about filtered strategies (bold, italic, accent color) ---
$ws.Range("E4").Value = "Filtered out skipped strategies"
$ws.Range("E4").Font.Bold = $true
$ws.Range("E4").Font.Italic = $true
$ws.Range("E4").Font.ThemeColor = 6

# --- Populate Discussion (Hypothesis), Potential Strategy and Strategy Group columns ---
$ws.Range("C7").Value = "Employees with high OT are leaving"
$ws.Range("D7").Value = "Reduce overtime"
$ws.Range("C8").Value = "Employees with Job level 1 are leaving and Job level 2 are staying"
$ws.Range("D8").Value = "Promote faster for high performers"
$ws.Range("C9").Value = "Low income - Bin 1 is likely to leave"
$ws.Range("D9").Value = "Skip - include in executive report"
$ws.Range("C10").Value = "Zero stock option - more likely to leave"
$ws.Range("D10").Value = "Skip - include in executive report"
$ws.Range("C11").Value = "If Years at company is high, they are more likely to stay. If low they are likely to leave"
$ws.Range("D11").Value = "Tie promotion if low to advance faster / Mentor if years at company low"
$ws.Range("E11").Value = "Professional development / Personal development"
$ws.Range("C12").Value = "Single employees are more likely to leave"
$ws.Range("D12").Value = "Skip - include in executive report"
$ws.Range("C13").Value = "If total working years is high, they are more likely to stay.If low they are likely to leave"
$ws.Range("D13").Value = "Tie low total woirking years to training and formation activities"
$ws.Range("E13").Value = "Personal development"
$ws.Range("C14").Value = "See total working years"
$ws.Range("D14").Value = "Skip - covered by total working years"
$ws.Range("C15").Value = "More time in current role related to lower attrition"
$ws.Range("D15").Value = "Incentivise specialisation or promote  / Mentorships role"
$ws.Range("E15").Value = "Professional development / Personal development"
$ws.Range("C16").Value = "Certain job roles have higher attrition - need to be monitored more closely"
$ws.Range("D16").Value = "Skip - include in executive report"
$ws.Range("C17").Value = "See years in current role"
$ws.Range("D17").Value = "Skip - covered by years in current role"
$ws.Range("C18").Value = "Employees with low environment satisfaction are more likely to leave"
$ws.Range("D18").Value = "Improve the workplace environment (address issues quoted ny those with value and low satisfaction)"
$ws.Range("E18").Value = "Work life"
$ws.Range("C19").Value = "Work life balance Bad - more likely to leave"
$ws.Range("D19").Value = "Improve worklife balance"
$ws.Range("E19").Value = "Work life"
$ws.Range("C20").Value = "More business travel - more likely to leave than if have less business travel"
$ws.Range("D20").Value = "Recduce business travel where possible"
$ws.Range("E20").Value = "Work life"
$ws.Range("C21").Value = "High Job involvement - likely to stay, If low they are likely to leave"
$ws.Range("D21").Value = "Create personal development plan"
$ws.Range("E21").Value = "Personal development"
$ws.Range("C22").Value = "If numbers of companies is very high (bin 4) more likely to leave than if above average (BIN 3)- out of our control"
$ws.Range("D22").Value = "Skip - include in executive report"
$ws.Range("C23").Value = "Low Job satisfaction - more likely to leave / High  Job satisfaction more likely to stay"
$ws.Range("D23").Value = "Low: create personal development plan, High: Suggest take on mentorship role"
$ws.Range("E23").Value = "Personal development"
$ws.Range("C24").Value = "Sales - more likely to leave / R&D mnore likely to stay"
$ws.Range("D24").Value = "Skip - include in executive report"
$ws.Range("C25").Value = "Educational field - technical degree slightly more likely to leave. More education might help to retain"
$ws.Range("D25").Value = "Skip - include in executive report"
$ws.Range("C26").Value = "Bin 1 - more likely to leave"
$ws.Range("D26").Value = "Skip - include in executive report"
$ws.Range("C27").Value = "High distance from home - more likely to leave"
$ws.Range("D27").Value = "Monitor worklife balance"
$ws.Range("E27").Value = "Work life"

# --- Apply the table AutoFilter on "Strategy Group" (4th column) to only show
#     Personal development related strategies; this also hides the other rows ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Range.AutoFilter(4, @("Personal development","Professional development / Personal development"), 7)

# --- Update active selection to reflect where the user left off ---
$ws.Range("D21").Select()
